# Auto-generated edit script: updates cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.590.63'
$ws.Range("E2").Value = '  -0.76%  '

# Row 3
$ws.Range("D3").Value = '2.294.53'
$ws.Range("E3").Value = '  -0.30%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.80'
$ws.Range("E5").Value = '  -1.68%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.87'
$ws.Range("E6").Value = '  -1.60%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").Value = '  -1.12%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  -2.00%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.65'
$ws.Range("E10").Value = '  -3.09%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.16'
$ws.Range("E11").Value = '  +5.44%  '

# Row 12
$ws.Range("E12").Value = '  -1.45%  '

# Row 13
$ws.Range("E13").Value = '  +0.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.78'
$ws.Range("E14").Value = '  +0.02%  '

# Row 15
$ws.Range("D15").Value = '2.651.84'
$ws.Range("E15").Value = '  -0.29%  '

# Row 16
$ws.Range("D16").Value = '2.295.79'
$ws.Range("E16").Value = '  -0.56%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.781'
$ws.Range("E17").Value = '  -0.37%  '

# Row 18
$ws.Range("D18").Value = '42.551.44'
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.17'
$ws.Range("E19").Value = '  -5.95%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  -1.86%  '

# Row 21
$ws.Range("E21").Value = '  -0.63%  '

# Row 22
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.81'
$ws.Range("E23").Value = '  -0.82%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("E24").Value = '  +2.92%  '

# Row 25
$ws.Range("E25").Value = '  +0.06%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  -2.59%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.62'
$ws.Range("E27").Value = '  -3.35%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.05'
$ws.Range("E28").Value = '  -0.05%  '

# Row 29
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '164.40'
$ws.Range("E29").Value = '  -1.44%  '

# Row 30
$ws.Range("E30").Value = '  -0.62%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.19'
$ws.Range("E31").Value = '  -2.48%  '

# Row 32
$ws.Range("E32").Value = '  -0.06%  '

# Row 33
$ws.Range("E33").Value = '  -1.29%  '

# Row 34
$ws.Range("E34").Value = '  -0.03%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.40'
$ws.Range("E35").Value = '  -8.53%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0695'
$ws.Range("E36").Value = '  +0.08%  '

# Row 37
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  -2.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0999'
$ws.Range("E38").Value = '  -1.65%  '

# Row 39
$ws.Range("E39").Value = '  -0.52%  '

# Row 40
$ws.Range("E40").Value = '  -0.25%  '

# Row 41
$ws.Range("E41").Value = '  -1.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.38'
$ws.Range("E42").Value = '  +7.76%  '

# Row 43
$ws.Range("D43").Value = '1.955.70'
$ws.Range("E43").Value = '  -2.89%  '

# Row 44
$ws.Range("E44").Value = '  +4.81%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0279'
$ws.Range("E45").Value = '  -1.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.05'
$ws.Range("E46").Value = '  -3.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.75'
$ws.Range("E47").Value = '  -1.46%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.91'
$ws.Range("E48").Value = '  -0.11%  '

# Row 49
$ws.Range("D49").Value = '2.520.48'
$ws.Range("E49").Value = '  -0.26%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.99'
$ws.Range("E50").Value = '  -2.03%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.56'
$ws.Range("E51").Value = '  -0.82%  '
